$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at T (column 20), shifting nom/url_produit (T,U) to (U,V)
$ws.Columns.Item(20).Insert()

# New header for the inserted column T1
$ws.Range("T1").Value = "2026-01-28 13:40:26"

# For each data row, the new column T gets a copy of the last price column (S)
# when S has a numeric value; otherwise it stays blank (matches source data).
$lastRow = 205
for ($r = 2; $r -le $lastRow; $r++) {
    $sVal = $ws.Cells.Item($r, 19).Value2
    if ($sVal -ne $null -and $sVal -ne "") {
        $ws.Cells.Item($r, 20).Value = $sVal
    }
}
